# Add two new columns, I ("I0") and J ("IF"), to the right of the existing
# data table (which currently ends at column H / "IP").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Write the new header labels.
$ws.Cells.Item(1, 9).Value  = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Match the existing header formatting (bold, centered, thin box border)
# by copying the format from the existing "IP" header (H1) onto the two
# new header cells, rather than re-building the style by hand.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-write the values after the paste (PasteSpecial of formats only should
# leave them untouched, but make sure the text is correct regardless of
# paste behavior).
$ws.Cells.Item(1, 9).Value  = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# --- Data rows (2-29) ---------------------------------------------------
$iVals = @(4,2,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,6,1,1,1,1,1,1,1,1)
$jVals = @(6,5,5,4,2,5,6,5,6,6,5,3,6,5,4,6,7,5,5,9,5,3,6,1,5,4,4,2)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value  = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
